# This script applies the "Update countries & provincias Spain" edit:
# - refreshes the "Datos actualizados" timestamp in A1
# - re-sorts/refreshes a subset of country rows (descending by Casos totales)
#   by overwriting the Pais/Casos totales/Nuevos casos/Casos activos/
#   Recuperados/Casos criticos/Muertes hoy/Muertes columns in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "Datos actualizados" timestamp cell
$ws.Range("A1").Value = 'Datos actualizados a 14 de Mayo de 2020 a las 15:05'

# 2) Rows that changed (country name and/or numeric columns) after the
#    data refresh / re-sort. Each entry: row, Pais, B..H
$rows = @(
    @(4, 'Estados Unidos', 1430653,305,310259,1035160,16349,37,85234),
    @(9, 'Brasil', 192081,2924,78424,100381,8318,118,13276),
    @(19, 'Arabia Saudita', 46869,2039,19051,27535,147,10,283),
    @(27, 'Portugal', 28319,187,3198,23937,108,9,1184),
    @(28, 'Catar', 28272,1733,3356,24902,72,0,14),
    @(49, 'Serbia', 10374,79,4084,6066,22,2,224),
    @(58, 'Finlandia', 6145,91,4300,1558,33,3,287),
    @(79, 'Bosnia y Herzegovina', 2218,37,1272,824,4,2,122),
    @(105, 'Tayikistan', 907,106,0,878,0,6,29),
    @(106, 'Republica de Chipre', 905,0,449,439,10,0,17),
    @(107, 'Albania', 898,18,694,173,1,0,31),
    @(108, 'Libano', 886,8,236,624,4,0,26),
    @(109, 'Niger', 860,0,658,153,0,0,49),
    @(110, 'Guinea-Bisau', 836,0,26,807,0,0,3),
    @(111, 'Costa Rica', 815,0,527,280,6,1,8),
    @(114, 'Kenia', 758,21,284,432,1,2,42),
    @(115, 'Mali', 758,0,412,302,0,0,44),
    @(116, 'Paraguay', 740,0,182,547,0,0,11),
    @(119, 'Zambia', 654,208,124,523,1,0,7),
    @(120, 'Georgia', 652,5,383,257,6,1,12),
    @(121, 'San Marino', 648,5,172,435,1,0,41),
    @(122, 'Jordania', 582,0,392,181,5,0,9),
    @(123, 'Guinea Ecuatorial', 522,0,13,503,0,0,6),
    @(124, 'Malta', 522,14,443,73,1,0,6),
    @(125, 'Jamaica', 509,2,113,387,0,0,9),
    @(126, 'Tanzania', 509,0,183,305,7,0,21),
    @(193, 'Nueva Caledonia', 18,0,18,0,0,0,0),
    @(195, 'Santa Lucia', 18,0,18,0,0,0,0),
    @(200, 'Dominica', 16,0,15,1,0,0,0),
    @(201, 'Curazao', 16,0,14,1,0,0,1),
    @(214, 'Bonaire, San Eustaquio y Saba', 6,0,6,0,0,0,0),
    @(215, 'San Bartolome', 6,0,6,0,0,0,0),
    @(216, 'Sahara Occidental', 6,0,6,0,0,0,0)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
}
